# Apply the edit described in the diff:
#  - Add a blank separator row (42)
#  - Add a second 15-row data block (rows 43-57), same layout/formatting
#    as the existing block in rows 3-17, but only populating columns B:M
#    (values 1..12 repeating) - columns N:P keep the formatting but no value.
#  - Add a second 15-row "calculated" block (rows 63-77) with a fixed
#    SUM formula (explicit neighbour refs instead of a merged range, per
#    the commit message "arreglo para b en caso de 1 proceso").
#  - Update the sheet view (scrolled/selected range) and the window state
#    (workbook minimized).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Blank separator row, mirroring row 2 ---------------------------
$ws.Rows(42).RowHeight = 15

# --- 2. Second data block (rows 43-57), copied formatting from 3-17 ----
$ws.Range("B3:P17").Copy()
$ws.Range("B43:P57").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

for ($r = 43; $r -le 57; $r++) {
    for ($c = 2; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = $c - 1
    }
}

# --- 3. Second calculated block (rows 63-77) ----------------------------
# Same relative pattern as the existing rows 23-37 block, but written with
# explicit neighbour references (this is the "fix" from the commit
# message) instead of a merged range, so it still works when there is
# only one data column.
$ws.Range("B63:M77").FormulaR1C1 = "=SUM(R[-20]C,R[-21]C,R[-20]C[-1],R[-19]C,R[-20]C[1])"

# --- 4. View state -------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("B63:M77").Select()

$excel.ActiveWindow.WindowState = -4140
